# Re-shuffle the species-record fields (Id, Taxonsorteringsordning,
# Rodlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Antal, Ost,
# Nord, Noggrannhet) across rows 3-16. All other fields (which are
# identical for every record in this sheet) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 111902031
$ws.Range("B3").Value = 90660
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 4362
$ws.Range("F3").Value = "Blå taggsvamp"
$ws.Range("G3").Value = "Hydnellum caeruleum"
$ws.Range("H3").Value = "(Hornem.) P.Karst."
$ws.Range("I3").Value = "2"
$ws.Range("Q3").Value = 524990.2026765908
$ws.Range("R3").Value = 6867385.898910107
$ws.Range("S3").Value = 25

# Row 4
$ws.Range("A4").Value = 111902033
$ws.Range("B4").Value = 90300
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4745
$ws.Range("F4").Value = "Tallriska"
$ws.Range("G4").Value = "Lactarius musteus"
$ws.Range("H4").Value = "Fr."
$ws.Range("I4").Value = "1"
$ws.Range("Q4").Value = 525027.0938798942
$ws.Range("R4").Value = 6867370.16309081
$ws.Range("S4").Value = 10

# Row 5
$ws.Range("A5").Value = 111902035
$ws.Range("B5").Value = 90658
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 4361
$ws.Range("F5").Value = "Orange taggsvamp"
$ws.Range("G5").Value = "Hydnellum aurantiacum"
$ws.Range("H5").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I5").Value = "3"
$ws.Range("Q5").Value = 525047.2558985724
$ws.Range("R5").Value = 6867385.376238698
$ws.Range("S5").Value = 25

# Row 6
$ws.Range("A6").Value = 111902029
$ws.Range("B6").Value = 88032
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 6276
$ws.Range("F6").Value = "Goliatmusseron"
$ws.Range("G6").Value = "Tricholoma matsutake"
$ws.Range("H6").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I6").Value = "4"
$ws.Range("Q6").Value = 524971.6686743505
$ws.Range("R6").Value = 6867341.509407703
$ws.Range("S6").Value = 5

# Row 7
$ws.Range("A7").Value = 111902036
$ws.Range("B7").Value = 88032
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 6276
$ws.Range("F7").Value = "Goliatmusseron"
$ws.Range("G7").Value = "Tricholoma matsutake"
$ws.Range("H7").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I7").Value = "2"
$ws.Range("Q7").Value = 525015.987664115
$ws.Range("R7").Value = 6867405.860822954
$ws.Range("S7").Value = 25

# Row 8
$ws.Range("A8").Value = 111902028
$ws.Range("B8").Value = 90666
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I8").Value = "1"
$ws.Range("Q8").Value = 524954.0254130038
$ws.Range("R8").Value = 6867304.187839299
$ws.Range("S8").Value = 5

# Row 9
$ws.Range("A9").Value = 111902038
$ws.Range("B9").Value = 90666
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4364
$ws.Range("F9").Value = "Dropptaggsvamp"
$ws.Range("G9").Value = "Hydnellum ferrugineum"
$ws.Range("H9").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I9").Value = "1"
$ws.Range("Q9").Value = 524892.725176702
$ws.Range("R9").Value = 6867498.641564975
$ws.Range("S9").Value = 10

# Row 10
$ws.Range("A10").Value = 111902030
$ws.Range("B10").Value = 88032
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 6276
$ws.Range("F10").Value = "Goliatmusseron"
$ws.Range("G10").Value = "Tricholoma matsutake"
$ws.Range("H10").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I10").Value = "6"
$ws.Range("Q10").Value = 524971.3961406752
$ws.Range("R10").Value = 6867378.699329315
$ws.Range("S10").Value = 5

# Row 11
$ws.Range("A11").Value = 111902039
$ws.Range("B11").Value = 90682
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 2059
$ws.Range("F11").Value = "Skrovlig taggsvamp"
$ws.Range("G11").Value = "Hydnellum scabrosum"
$ws.Range("H11").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I11").Value = "5"
$ws.Range("Q11").Value = 524868.0170565489
$ws.Range("R11").Value = 6867460.329015278
$ws.Range("S11").Value = 5

# Row 12
$ws.Range("A12").Value = 111902037
$ws.Range("B12").Value = 90654
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 149
$ws.Range("F12").Value = "Tallgråticka"
$ws.Range("G12").Value = "Boletopsis grisea"
$ws.Range("H12").Value = "(Peck) Bondartsev & Singer"
$ws.Range("I12").Value = "2"
$ws.Range("Q12").Value = 524868.6293626219
$ws.Range("R12").Value = 6867441.031870116
$ws.Range("S12").Value = 5

# Row 13
$ws.Range("A13").Value = 111902032
$ws.Range("B13").Value = 90658
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 4361
$ws.Range("F13").Value = "Orange taggsvamp"
$ws.Range("G13").Value = "Hydnellum aurantiacum"
$ws.Range("H13").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I13").Value = "1"
$ws.Range("Q13").Value = 524989.2701192262
$ws.Range("R13").Value = 6867384.479730026
$ws.Range("S13").Value = 5

# Row 14
$ws.Range("A14").Value = 111902034
$ws.Range("B14").Value = 90660
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 4362
$ws.Range("F14").Value = "Blå taggsvamp"
$ws.Range("G14").Value = "Hydnellum caeruleum"
$ws.Range("H14").Value = "(Hornem.) P.Karst."
$ws.Range("I14").Value = "10"
$ws.Range("Q14").Value = 525038.6070930503
$ws.Range("R14").Value = 6867407.439287313
$ws.Range("S14").Value = 25

# Row 15
$ws.Range("A15").Value = 111902027
$ws.Range("B15").Value = 90660
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 4362
$ws.Range("F15").Value = "Blå taggsvamp"
$ws.Range("G15").Value = "Hydnellum caeruleum"
$ws.Range("H15").Value = "(Hornem.) P.Karst."
$ws.Range("I15").Value = "5"
$ws.Range("Q15").Value = 524936.9216418237
$ws.Range("R15").Value = 6867321.952660743
$ws.Range("S15").Value = 25

# Row 16
$ws.Range("A16").Value = 111902026
$ws.Range("B16").Value = 90682
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 2059
$ws.Range("F16").Value = "Skrovlig taggsvamp"
$ws.Range("G16").Value = "Hydnellum scabrosum"
$ws.Range("H16").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q16").Value = 524951.0483835863
$ws.Range("R16").Value = 6867324.410012136
$ws.Range("S16").Value = 10
